$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.206.76"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.681.71"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5248"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2692"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07632"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.22%  "

$ws.Range("D12").Value = "1.677.84"
$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.515"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5747"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008255"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.30%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.07"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("D17").Value = "26.246.76"
$ws.Range("E17").Value = "  -0.13%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.867"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.62%  "

$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.239"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "7.787"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1260"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.83"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.61%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06299"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.315"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.571"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.42%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.682"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.024"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6114"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.422"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.79%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.750"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.10%  "

$ws.Range("E38").Value = "  -0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.56%  "

$ws.Range("D40").Value = "1.097.87"
$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8841"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.38%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.41"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.06%  "

$ws.Range("D44").Value = "1.831.13"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.005"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.085"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05273"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4279"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.012"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.03%  "
